# Updates the "Program Control" and "Variables" worksheets so that the
# "GoTo" / "Comments" (line / var type) columns are now populated, and the
# related summary counts / averages on the "Variables" sheet are refreshed.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Program Control" -- a handful of For/If counts were corrected
# ----------------------------------------------------------------------
$wsControl = $wb.Worksheets.Item("Program Control")

$wsControl.Cells.Item(9, 3).Value  = 4   # C9
$wsControl.Cells.Item(9, 5).Value  = 3   # E9

$wsControl.Cells.Item(10, 3).Value = 4   # C10
$wsControl.Cells.Item(10, 5).Value = 3   # E10

$wsControl.Cells.Item(11, 3).Value = 1   # C11
$wsControl.Cells.Item(11, 5).Value = 1   # E11

$wsControl.Cells.Item(13, 3).Value = 2   # C13

$wsControl.Cells.Item(14, 3).Value = 4   # C14
$wsControl.Cells.Item(14, 5).Value = 3   # E14

$wsControl.Cells.Item(21, 3).Value = 4   # C21
$wsControl.Cells.Item(21, 5).Value = 3   # E21

$wsControl.Cells.Item(22, 5).Value = 7   # E22

$wsControl.Cells.Item(23, 3).Value = 2   # C23

$wsControl.Cells.Item(24, 3).Value = 4   # C24
$wsControl.Cells.Item(24, 5).Value = 6   # E24
$wsControl.Cells.Item(24, 6).Value = 2   # F24

$wsControl.Cells.Item(26, 5).Value = 3   # E26

$wsControl.Cells.Item(27, 3).Value = 1   # C27
$wsControl.Cells.Item(27, 5).Value = 1   # E27

$wsControl.Cells.Item(28, 2).Value = 2   # B28
$wsControl.Cells.Item(28, 3).Value = 6   # C28
$wsControl.Cells.Item(28, 5).Value = 5   # E28

# ----------------------------------------------------------------------
# Sheet 2: "Variables" -- the "Comments" (var type) column H is now
# printed, and B/C counts + the average length (G) were recomputed.
# ----------------------------------------------------------------------
$wsVars = $wb.Worksheets.Item("Variables")

$wsVars.Cells.Item(2, 8).Value  = 1                    # H2

$wsVars.Cells.Item(3, 8).Value  = 1                    # H3

$wsVars.Cells.Item(4, 8).Value  = 2                    # H4

$wsVars.Cells.Item(5, 3).Value  = 8                    # C5
$wsVars.Cells.Item(5, 5).Value  = 4                    # E5
$wsVars.Cells.Item(5, 7).Value  = 2.875                # G5
$wsVars.Cells.Item(5, 8).Value  = 0                    # H5

$wsVars.Cells.Item(6, 7).Value  = 2.333333333333333    # G6
$wsVars.Cells.Item(6, 8).Value  = 1                    # H6

$wsVars.Cells.Item(7, 7).Value  = 9.090909090909092    # G7
$wsVars.Cells.Item(7, 8).Value  = 1                    # H7

$wsVars.Cells.Item(8, 7).Value  = 8.642857142857142    # G8
$wsVars.Cells.Item(8, 8).Value  = 0                    # H8

$wsVars.Cells.Item(9, 7).Value  = 2.75                 # G9
$wsVars.Cells.Item(9, 8).Value  = 1                    # H9

$wsVars.Cells.Item(10, 7).Value = 3                    # G10
$wsVars.Cells.Item(10, 8).Value = 0                    # H10

$wsVars.Cells.Item(11, 7).Value = 1                    # G11
$wsVars.Cells.Item(11, 8).Value = 1                    # H11

$wsVars.Cells.Item(12, 7).Value = 3.888888888888889    # G12
$wsVars.Cells.Item(12, 8).Value = 0                    # H12

$wsVars.Cells.Item(13, 7).Value = 2.333333333333333    # G13
$wsVars.Cells.Item(13, 8).Value = 0                    # H13

$wsVars.Cells.Item(14, 8).Value = 0                    # H14

$wsVars.Cells.Item(15, 7).Value = 4.142857142857143    # G15
$wsVars.Cells.Item(15, 8).Value = 2                    # H15

$wsVars.Cells.Item(16, 2).Value = 15                   # B16
$wsVars.Cells.Item(16, 3).Value = 6                    # C16
$wsVars.Cells.Item(16, 7).Value = 6.266666666666667    # G16
$wsVars.Cells.Item(16, 8).Value = 0                    # H16

$wsVars.Cells.Item(17, 7).Value = 6                    # G17
$wsVars.Cells.Item(17, 8).Value = 0                    # H17

$wsVars.Cells.Item(18, 2).Value = 17                   # B18
$wsVars.Cells.Item(18, 3).Value = 8                    # C18
$wsVars.Cells.Item(18, 7).Value = 5.823529411764706    # G18
$wsVars.Cells.Item(18, 8).Value = 0                    # H18

$wsVars.Cells.Item(19, 7).Value = 6.588235294117647    # G19
$wsVars.Cells.Item(19, 8).Value = 0                    # H19

$wsVars.Cells.Item(20, 2).Value = 18                   # B20
$wsVars.Cells.Item(20, 3).Value = 8                    # C20
$wsVars.Cells.Item(20, 7).Value = 3.888888888888889    # G20
$wsVars.Cells.Item(20, 8).Value = 0                    # H20

$wsVars.Cells.Item(21, 7).Value = 5.769230769230769    # G21
$wsVars.Cells.Item(21, 8).Value = 0                    # H21

$wsVars.Cells.Item(22, 2).Value = 18                   # B22
$wsVars.Cells.Item(22, 3).Value = 10                   # C22
$wsVars.Cells.Item(22, 7).Value = 3.833333333333333    # G22
$wsVars.Cells.Item(22, 8).Value = 2                    # H22

$wsVars.Cells.Item(23, 2).Value = 8                    # B23
$wsVars.Cells.Item(23, 3).Value = 4                    # C23
$wsVars.Cells.Item(23, 5).Value = 3                    # E23
$wsVars.Cells.Item(23, 7).Value = 2.5                  # G23
$wsVars.Cells.Item(23, 8).Value = 0                    # H23

$wsVars.Cells.Item(24, 2).Value = 16                   # B24
$wsVars.Cells.Item(24, 3).Value = 3                    # C24
$wsVars.Cells.Item(24, 5).Value = 7                    # E24
$wsVars.Cells.Item(24, 6).Value = 6                    # F24
$wsVars.Cells.Item(24, 7).Value = 5.5                  # G24
$wsVars.Cells.Item(24, 8).Value = 1                    # H24

$wsVars.Cells.Item(25, 7).Value = 4.833333333333333    # G25
$wsVars.Cells.Item(25, 8).Value = 1                    # H25

$wsVars.Cells.Item(26, 2).Value = 16                   # B26
$wsVars.Cells.Item(26, 3).Value = 8                    # C26
$wsVars.Cells.Item(26, 7).Value = 2.0625               # G26
$wsVars.Cells.Item(26, 8).Value = 0                    # H26

$wsVars.Cells.Item(27, 7).Value = 3.4                  # G27
$wsVars.Cells.Item(27, 8).Value = 0                    # H27

$wsVars.Cells.Item(28, 2).Value = 22                   # B28
$wsVars.Cells.Item(28, 3).Value = 9                    # C28
$wsVars.Cells.Item(28, 4).Value = 0                    # D28
$wsVars.Cells.Item(28, 5).Value = 6                    # E28
$wsVars.Cells.Item(28, 6).Value = 7                    # F28
$wsVars.Cells.Item(28, 7).Value = 5.636363636363637    # G28
$wsVars.Cells.Item(28, 8).Value = 1                    # H28
